$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C28").Value = 461
$ws.Range("D28").Value = 50
$ws.Range("E28").Value = 411
$ws.Range("F28").Value = 7.78816199376947
